$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: replace sample data with new sample data
$ws.Range("A2").Value = 20005890
$ws.Range("B2").Value = "'2017-12-01"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "AA"
$ws.Range("F2").Value = 0

# Rows 3-6: clear all values (A, E, F fully empty; B, C, D keep formatting but no value)
$ws.Range("A3:F6").ClearContents()

# Set the active selection to E3
$ws.Range("E3").Select()
